$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.422.56"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.861.07"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4743"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2748"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06446"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "1.848.52"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07443"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.013"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6372"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.82%  "
$ws.Range("D16").Value = "30.365.06"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.34%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007420"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").Value = "2.099.94"
$ws.Range("E21").Value = "  -4.01%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.024"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.022"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.280"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("E29").Value = "  +7.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.391"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.157"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.944"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04919"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.156"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7296"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9994"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.694"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01899"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.652"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9122"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.978"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4126"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.589"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.154"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1215"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.738"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.414"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.46%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.50%  "
